# Update gh-pages to output generated at 456a3b4
# Refresh the "views" counters (column F) across all four sheets so the
# numbers reflect the newly generated output.

$wb = $excel.ActiveWorkbook

function Set-Views {
    param(
        [string]$SheetName,
        [hashtable]$Updates
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $Updates.Keys) {
        $ws.Range($cellRef).Value = $Updates[$cellRef]
    }
}

Set-Views "展览" @{
    "F3"  = 44
    "F4"  = 5922
    "F13" = 1591
    "F14" = 1591
    "F15" = 1580
    "F16" = 555
    "F17" = 167
    "F18" = 628
    "F19" = 4500
    "F20" = 33
    "F22" = 3346
    "F23" = 818
    "F24" = 20
    "F27" = 2317
    "F29" = 339
    "F31" = 457
    "F32" = 1233
    "F34" = 10
    "F35" = 5
    "F37" = 1222
    "F38" = 1201
}

Set-Views "演出" @{
    "F15" = 68
    "F18" = 127
    "F19" = 304
    "F20" = 233
    "F21" = 495
}

Set-Views "本地生活" @{
    "F2" = 250
    "F3" = 698
    "F4" = 190
    "F5" = 279
}

Set-Views "全部类型" @{
    "F2"  = 250
    "F4"  = 44
    "F5"  = 698
    "F6"  = 190
    "F7"  = 5922
    "F8"  = 5922
    "F23" = 1591
    "F25" = 1580
    "F26" = 555
    "F27" = 167
    "F28" = 628
    "F29" = 4500
    "F31" = 3346
    "F32" = 818
    "F33" = 20
    "F36" = 2317
    "F38" = 339
    "F40" = 457
    "F41" = 1233
    "F42" = 127
    "F43" = 304
    "F44" = 233
    "F45" = 495
    "F48" = 1222
    "F50" = 1201
}
